$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns per latest scrape
$ws.Range("D2").Value = "'27.333.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.15%  "
$ws.Range("D3").Value = "'1.856.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.99%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'323.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.4522"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.22%  "
$ws.Range("E8").Value = "  -4.76%  "
$ws.Range("D9").Value = "'48.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.22%  "
$ws.Range("D10").Value = "'0.07924"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.35%  "
$ws.Range("D11").Value = "'1.014"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.60%  "
$ws.Range("D12").Value = "'21.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("D13").Value = "'1.868.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.23%  "
$ws.Range("D14").Value = "'5.902"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.64%  "
$ws.Range("D15").Value = "'7.121"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.70%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "'85.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.75%  "
$ws.Range("D18").Value = "'0.00001027"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.79%  "
$ws.Range("D19").Value = "'0.06542"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "'17.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.69%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'5.536"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.37%  "
$ws.Range("D23").Value = "'27.333.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("D24").Value = "'10.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.88%  "
$ws.Range("D25").Value = "'2.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "'2.082.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("D27").Value = "'153.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "'2.064"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").Value = "'5.427"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.10%  "
$ws.Range("D31").Value = "'121.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").Value = "'1.478"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("D33").Value = "'0.09280"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.54%  "
$ws.Range("D34").Value = "'0.9330"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("D35").Value = "'3.601"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").Value = "'5.255"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.01%  "
$ws.Range("D37").Value = "'1.226"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").Value = "'0.02224"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("D39").Value = "'0.05980"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'8.138"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.58%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "'0.5900"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.74%  "
$ws.Range("D43").Value = "'0.1891"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "'10.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.96%  "
$ws.Range("D45").Value = "'1.284"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").Value = "'0.5599"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.27%  "
$ws.Range("D47").Value = "'11.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.85%  "
$ws.Range("D48").Value = "'3.375"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D49").Value = "'1.913"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.59%  "
$ws.Range("D50").Value = "'0.06755"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "'108.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.12%  "
